$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row (row 1) cells: "_old" suffix -> "_FV2210", "_new" suffix -> "_FV2304"
$leftHeaders  = @("Segmentname_FV2210","Segmentgruppe_FV2210","Segment_FV2210","Datenelement_FV2210","Segment ID_FV2210","Code_FV2210","Qualifier_FV2210","Beschreibung_FV2210","Bedingungsausdruck_FV2210","Bedingung_FV2210")
$rightHeaders = @("Segmentname_FV2304","Segmentgruppe_FV2304","Segment_FV2304","Datenelement_FV2304","Segment ID_FV2304","Code_FV2304","Qualifier_FV2304","Beschreibung_FV2304","Bedingungsausdruck_FV2304","Bedingung_FV2304")

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $leftHeaders[$i]
}
# Column 11 (K1) stays "diff" - unchanged
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $rightHeaders[$i]
}

# 2. Turn the A1:U64 range into an Excel Table ("Table1") with an AutoFilter
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U64"), $null, 1)
$tbl.Name = "Table1"

# 3. Freeze the header row (split after row 1)
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
